$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value2 = "Volume 30   Number  6"
$ws.Range("C9").Value2 = "Report Covering the Week  2/6/2023  Through  2/12/2023"

# --- Data table updates ---
# Row 15
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("C15").Value2 = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value2 = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value2 = 0
$ws.Range("F15").NumberFormat = '#,##0'
$ws.Range("F15").Value2 = 1
$ws.Range("G15").Value2 = 3
$ws.Range("H15").Value2 = -66.666666666666
$ws.Range("I15").NumberFormat = '#,##0'
$ws.Range("I15").Value2 = 1
$ws.Range("J15").Value2 = 3
$ws.Range("K15").Value2 = -66.666666666666
$ws.Range("L15").Value2 = 0
$ws.Range("M15").Value2 = -50
$ws.Range("N15").Value2 = -66.666666666666
# Row 16
$ws.Range("C16").Value2 = 2
$ws.Range("D16").Value2 = 4
$ws.Range("E16").Value2 = -50
$ws.Range("F16").Value2 = 13
$ws.Range("G16").Value2 = 14
$ws.Range("H16").Value2 = -7.142857142857
$ws.Range("I16").Value2 = 17
$ws.Range("J16").Value2 = 21
$ws.Range("K16").Value2 = -19.047619047619
$ws.Range("L16").Value2 = 41.666666666666
$ws.Range("M16").Value2 = -55.263157894736
$ws.Range("N16").Value2 = -86.507936507936
# Row 17
$ws.Range("C17").Value2 = 10
$ws.Range("D17").Value2 = 7
$ws.Range("E17").Value2 = 42.857142857142
$ws.Range("F17").Value2 = 31
$ws.Range("G17").Value2 = 30
$ws.Range("H17").Value2 = 3.333333333333
$ws.Range("I17").Value2 = 50
$ws.Range("J17").Value2 = 49
$ws.Range("K17").Value2 = 2.040816326530
$ws.Range("L17").Value2 = 28.205128205128
$ws.Range("M17").Value2 = 72.413793103448
$ws.Range("N17").Value2 = 21.951219512195
# Row 18
$ws.Range("C18").Value2 = 2
$ws.Range("D18").Value2 = 6
$ws.Range("E18").Value2 = -66.666666666666
$ws.Range("F18").Value2 = 5
$ws.Range("G18").Value2 = 18
$ws.Range("H18").Value2 = -72.222222222222
$ws.Range("I18").Value2 = 17
$ws.Range("J18").Value2 = 25
$ws.Range("K18").Value2 = -32
$ws.Range("L18").Value2 = -10.526315789473
$ws.Range("M18").Value2 = -65.306122448979
$ws.Range("N18").Value2 = -88.815789473684
# Row 19
$ws.Range("C19").Value2 = 25
$ws.Range("D19").Value2 = 11
$ws.Range("E19").Value2 = 127.272727272727
$ws.Range("F19").Value2 = 52
$ws.Range("G19").Value2 = 39
$ws.Range("H19").Value2 = 33.333333333333
$ws.Range("I19").Value2 = 68
$ws.Range("J19").Value2 = 57
$ws.Range("K19").Value2 = 19.298245614035
$ws.Range("L19").Value2 = 54.545454545454
$ws.Range("M19").Value2 = 38.775510204081
$ws.Range("N19").Value2 = 9.677419354838
# Row 20
$ws.Range("C20").Value2 = 6
$ws.Range("D20").Value2 = 7
$ws.Range("E20").Value2 = -14.285714285714
$ws.Range("F20").Value2 = 18
$ws.Range("G20").Value2 = 17
$ws.Range("H20").Value2 = 5.882352941176
$ws.Range("I20").Value2 = 26
$ws.Range("J20").Value2 = 21
$ws.Range("K20").Value2 = 23.809523809523
$ws.Range("L20").Value2 = 13.043478260869
$ws.Range("M20").Value2 = -46.938775510204
$ws.Range("N20").Value2 = -93.450881612090
# Row 21
$ws.Range("C21").Value2 = 46
$ws.Range("D21").Value2 = 36
$ws.Range("E21").Value2 = 27.777777777777
$ws.Range("F21").Value2 = 120
$ws.Range("G21").Value2 = 121
$ws.Range("H21").Value2 = -0.826446280991
$ws.Range("I21").Value2 = 179
$ws.Range("J21").Value2 = 176
$ws.Range("K21").Value2 = 1.704545454545
$ws.Range("L21").Value2 = 28.776978417266
$ws.Range("M21").Value2 = -17.129629629629
$ws.Range("N21").Value2 = -77.109974424552
# Row 24
$ws.Range("C24").Value2 = 30
$ws.Range("D24").Value2 = 20
$ws.Range("E24").Value2 = 50
$ws.Range("F24").Value2 = 96
$ws.Range("G24").Value2 = 93
$ws.Range("H24").Value2 = 3.225806451612
$ws.Range("I24").Value2 = 137
$ws.Range("J24").Value2 = 147
$ws.Range("K24").Value2 = -6.802721088435
$ws.Range("L24").Value2 = -8.666666666666
$ws.Range("M24").Value2 = 38.383838383838
# Row 25
$ws.Range("C25").Value2 = 6
$ws.Range("D25").Value2 = 5
$ws.Range("E25").Value2 = 20
$ws.Range("F25").Value2 = 36
$ws.Range("G25").Value2 = 39
$ws.Range("H25").Value2 = -7.692307692307
$ws.Range("I25").Value2 = 58
$ws.Range("J25").Value2 = 56
$ws.Range("K25").Value2 = 3.571428571428
$ws.Range("L25").Value2 = 81.25
$ws.Range("M25").Value2 = -21.621621621621
# Row 26
$ws.Range("C26").Value2 = 1
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("D26").Value2 = 2
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E26").Value2 = -50
$ws.Range("F26").Value2 = 2
$ws.Range("G26").Value2 = 4
$ws.Range("H26").Value2 = -50
$ws.Range("I26").Value2 = 3
$ws.Range("J26").Value2 = 4
$ws.Range("K26").Value2 = -25
$ws.Range("L26").Value2 = 50
# Row 27
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value2 = "0"
$ws.Range("D27").Value2 = 3
$ws.Range("E27").Value2 = -100
$ws.Range("F27").Value2 = 4
$ws.Range("G27").Value2 = 5
$ws.Range("H27").Value2 = -20
$ws.Range("I27").Value2 = 5
$ws.Range("J27").Value2 = 5
$ws.Range("K27").Value2 = 0
$ws.Range("L27").Value2 = 66.666666666666
# Row 28
$ws.Range("M28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M28").Value2 = 100
# Row 29
$ws.Range("M29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M29").Value2 = 0
